$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for every data row (2-514)
# from serial 46060 (2026-02-07) to serial 46061 (2026-02-08).
$ws.Range("C2:C514").Value = 46061
